$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(77).Insert()
